# Updated cryptos list on Sat Oct 12 20:42:59 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.057.03"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "2.470.92"
$ws.Range("E3").Value = "  +2.15%  "

$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").Value = "577.63"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").Value = "146.69"
$ws.Range("E6").Value = "  +0.50%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("E8").Value = "  -0.64%  "

$ws.Range("D9").Value = "2.471.56"
$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("D10").Value = "0.111"
$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("E11").Value = "  +1.68%  "

$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").Value = "28.82"
$ws.Range("E14").Value = "  +5.28%  "

$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("D16").Value = "2.920.27"
$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("D17").Value = "63.015.93"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").Value = "2.474.42"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("D19").Value = "8.19"
$ws.Range("E19").Value = "  +3.78%  "

$ws.Range("D20").Value = "11.05"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").Value = "329.38"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("E22").Value = "  +9.77%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").Value = "66.28"
$ws.Range("E25").Value = "  +0.89%  "

$ws.Range("D26").Value = "669.68"
$ws.Range("E26").Value = "  +4.64%  "

$ws.Range("D27").Value = "9.64"
$ws.Range("E27").Value = "  +13.26%  "

$ws.Range("D28").Value = "0.0₃0990"
$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("D29").Value = "2.593.16"
$ws.Range("E29").Value = "  +2.38%  "

$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -9.48%  "

$ws.Range("E31").Value = "  +2.61%  "

$ws.Range("D32").Value = "8.04"
$ws.Range("E32").Value = "  -1.92%  "

$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("E34").Value = "  -3.29%  "

$ws.Range("E35").Value = "  +3.60%  "

$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("E37").Value = "  +0.47%  "

$ws.Range("D38").Value = "5.45"
$ws.Range("E38").Value = "  +0.80%  "

$ws.Range("E39").Value = "  -0.69%  "

$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "18.77"
$ws.Range("E40").Value = "  +0.44%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "151.11"
$ws.Range("E41").Value = "  -1.21%  "

$ws.Range("D42").Value = "2.73"
$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("D43").Value = "1.75"
$ws.Range("E43").Value = "  -0.84%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").Value = "0.0₆0304"
$ws.Range("E45").Value = "  +5.88%  "

$ws.Range("D46").Value = "151.84"
$ws.Range("E46").Value = "  +4.79%  "

$ws.Range("D47").Value = "15.13"
$ws.Range("E47").Value = "  +20.88%  "

$ws.Range("D48").Value = "3.59"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("D49").Value = "'20.60"
$ws.Range("E49").Value = "  +0.83%  "

$ws.Range("E50").Value = "  +0.43%  "

$ws.Range("E51").Value = "  -0.93%  "
